# Updated cryptos list values (Price and Volume(1h) columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.079.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.646.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5217"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2607"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06319"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07683"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.645.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.421"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.871.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5575"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8237"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.096.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.732"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.214"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.444"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1212"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.399"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05894"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.268"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.419"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.409"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.655"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9880"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.395"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.760"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5665"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01620"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8597"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.838"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.030.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.798.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.107"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05185"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4225"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.32%  "
